$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.370.78'
$ws.Range('E2').Value = '  -3.14%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.250.47'
$ws.Range('E3').Value = '  -3.73%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '233.97'
$ws.Range('E5').Value = '  -2.04%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.634'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '69.82'
$ws.Range('E7').Value = '  -2.69%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.563'
$ws.Range('E9').Value = '  -3.70%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0993'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '58.38'
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '36.04'
$ws.Range('E12').Value = '  +11.85%  '
$ws.Range('E13').Value = '  -1.85%  '
$ws.Range('E14').Value = '  -4.70%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.584.97'
$ws.Range('E15').Value = '  -3.75%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.14'
$ws.Range('E16').Value = '  -6.06%  '
$ws.Range('E17').Value = '  -3.81%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.250.57'
$ws.Range('E18').Value = '  -3.97%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '42.277.83'
$ws.Range('E19').Value = '  -3.04%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0977'
$ws.Range('E20').Value = '  -3.05%  '
$ws.Range('E21').Value = '  -5.47%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '73.50'
$ws.Range('E22').Value = '  -5.97%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.68'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.00'
$ws.Range('E24').Value = '  +5.34%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.40'
$ws.Range('E27').Value = '  -3.32%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.06'
$ws.Range('E28').Value = '  -2.42%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  -3.74%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '169.26'
$ws.Range('E30').Value = '  -3.29%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '20.62'
$ws.Range('E31').Value = '  -6.69%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.121'
$ws.Range('E32').Value = '  -4.13%  '
$ws.Range('E33').Value = '  -5.49%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.40'
$ws.Range('E34').Value = '  +1.16%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0722'
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('E36').Value = '  -6.69%  '
$ws.Range('E37').Value = '  -2.73%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '21.69'
$ws.Range('E39').Value = '  -3.85%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.98'
$ws.Range('E40').Value = '  -5.58%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0269'
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '66.35'
$ws.Range('E42').Value = '  +2.97%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.93'
$ws.Range('E43').Value = '  -6.33%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.98'
$ws.Range('E44').Value = '  -2.41%  '
$ws.Range('E45').Value = '  -3.83%  '
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0₃0155'
$ws.Range('E48').Value = '  +22.36%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.44'
$ws.Range('E49').Value = '  +10.37%  '
$ws.Range('B50').Value = 'Celestia'
$ws.Range('C50').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '10.17'
$ws.Range('E50').Value = '  +9.63%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.18'
$ws.Range('E51').Value = '  -2.95%  '
